$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.837.81"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3
$ws.Range("D3").Value = "2.236.71"
$ws.Range("E3").Value = "  -1.66%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.10"
$ws.Range("E5").Value = "  +2.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.36"
$ws.Range("E6").Value = "  +4.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -1.61%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.24%  "

# Row 9
$ws.Range("E9").Value = "  +0.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.45"
$ws.Range("E10").Value = "  -0.89%  "

# Row 11
$ws.Range("E11").Value = "  -0.71%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.04"
$ws.Range("E12").Value = "  -0.94%  "

# Row 13
$ws.Range("E13").Value = "  -3.00%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.25"
$ws.Range("E14").Value = "  -0.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.873"
$ws.Range("E15").Value = "  +1.32%  "

# Row 16
$ws.Range("D16").Value = "2.576.39"
$ws.Range("E16").Value = "  -1.54%  "

# Row 17
$ws.Range("D17").Value = "2.235.71"
$ws.Range("E17").Value = "  -1.70%  "

# Row 18
$ws.Range("D18").Value = "42.833.84"
$ws.Range("E18").Value = "  -0.77%  "

# Row 19
$ws.Range("E19").Value = "  -0.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.76"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.16"
$ws.Range("E21").Value = "  +0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").Value = "  -3.82%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.99"
$ws.Range("E23").Value = "  +5.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "231.37"
$ws.Range("E24").Value = "  -1.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.24"
$ws.Range("E25").Value = "  -1.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.97"
$ws.Range("E26").Value = "  +5.94%  "

# Row 27
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.95"
$ws.Range("E28").Value = "  -0.74%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.28"
$ws.Range("E29").Value = "  -1.12%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.28"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.55"
$ws.Range("E32").Value = "  +0.14%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.06"
$ws.Range("E33").Value = "  -1.75%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0890"
$ws.Range("E34").Value = "  -0.86%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.57"
$ws.Range("E35").Value = "  -1.15%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.43"
$ws.Range("E36").Value = "  +13.80%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.128"
$ws.Range("E37").Value = "  -1.27%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0372"
$ws.Range("E38").Value = "  +0.20%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.65"
$ws.Range("E39").Value = "  +0.68%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").Value = "  +0.99%  "

# Row 41
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -1.29%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.00"
$ws.Range("E42").Value = "  -6.68%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  -1.69%  "

# Row 44
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").Value = "  -7.72%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.33"
$ws.Range("E46").Value = "  -2.61%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.65"
$ws.Range("E47").Value = "  -6.91%  "

# Row 48
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +1.12%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.44"
$ws.Range("E49").Value = "  -0.86%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("E50").Value = "  -0.61%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.35"
$ws.Range("E51").Value = "  -0.86%  "
